$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows("47:53").Insert()
Write-Output "done"
